# "Fruta / hortaliza, semanal"
# Updates the weekly Fecha (D) and Volumen/Precio (J,K,L,M,P) figures for
# rows 3-18 of the Espárragos sheet to the new reported values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 44503
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 900
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = 950
$ws.Range("P3").Value = 950

$ws.Range("D4").Value = 44523

$ws.Range("D5").Value = 44524
$ws.Range("K5").Value = 800
$ws.Range("L5").Value = 900
$ws.Range("M5").Value = 850
$ws.Range("P5").Value = 850

$ws.Range("D6").Value = 44517
$ws.Range("J6").Value = 500

$ws.Range("D7").Value = 44530
$ws.Range("J7").Value = 300

$ws.Range("D8").Value = 44512
$ws.Range("J8").Value = 600

$ws.Range("D9").Value = 44537
$ws.Range("J9").Value = 400

$ws.Range("D10").Value = 44505
$ws.Range("J10").Value = 440
$ws.Range("K10").Value = 900
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = 950
$ws.Range("P10").Value = 950

$ws.Range("D11").Value = 44525
$ws.Range("J11").Value = 360
$ws.Range("K11").Value = 800
$ws.Range("L11").Value = 900
$ws.Range("M11").Value = 850
$ws.Range("P11").Value = 850

$ws.Range("D12").Value = 44516
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 900
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = 950
$ws.Range("P12").Value = 950

$ws.Range("D13").Value = 44508

$ws.Range("D14").Value = 44532
$ws.Range("J14").Value = 240
$ws.Range("K14").Value = 800
$ws.Range("L14").Value = 900
$ws.Range("M14").Value = 850
$ws.Range("P14").Value = 850

$ws.Range("D15").Value = 44476
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 1100
$ws.Range("L15").Value = 1200
$ws.Range("M15").Value = 1150
$ws.Range("P15").Value = 1150

$ws.Range("D16").Value = 44510
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 900
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = 950
$ws.Range("P16").Value = 950

$ws.Range("D17").Value = 44511
$ws.Range("J17").Value = 500

$ws.Range("D18").Value = 44518
$ws.Range("J18").Value = 400
$ws.Range("K18").Value = 800
$ws.Range("L18").Value = 900
$ws.Range("M18").Value = 850
$ws.Range("P18").Value = 850
